$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="utf-8"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>TaittirIya BrAhmaNam TB 1.5 to 1.8  Book</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>Tamil</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>Corrections –</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t xml:space="preserve">Observed till </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:highlight w:val="red"/>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>????</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve">(ignore those which are already incorporated in your book’s version and date). </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:br/>
    <w:t>Kindly refer to your Guru for the differences in Swaram marking between various sources.</w:t>
  </w:r>
</w:p>
<w:tbl>
  <w:tblPr>
    <w:tblW w:w="13268" w:type="dxa"/>
    <w:tblInd w:w="-792" w:type="dxa"/>
    <w:tblBorders>
      <w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/>
      <w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/>
      <w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/>
      <w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/>
      <w:insideH w:val="single" w:sz="4" w:space="0" w:color="000000"/>
      <w:insideV w:val="single" w:sz="4" w:space="0" w:color="000000"/>
    </w:tblBorders>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="3310"/>
    <w:gridCol w:w="4738"/>
    <w:gridCol w:w="5220"/>
  </w:tblGrid>
  <w:tr>
  <w:tc>
    <w:tcPr>
      <w:tcW w:w="3310" w:type="dxa"/>
    </w:tcPr>
    <w:p>
      <w:pPr>
        <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>Section, Paragraph</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>Reference</w:t>
      </w:r>
    </w:p>
  </w:tc>
  <w:tc>
    <w:tcPr>
      <w:tcW w:w="4738" w:type="dxa"/>
    </w:tcPr>
    <w:p>
      <w:pPr>
        <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>As Printed</w:t>
      </w:r>
    </w:p>
  </w:tc>
  <w:tc>
    <w:tcPr>
      <w:tcW w:w="5220" w:type="dxa"/>
    </w:tcPr>
    <w:p>
      <w:pPr>
        <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:right="-18"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>To be read as or corrected as</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
        <w:ind w:right="-18"/>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>(Comments and analysis)</w:t>
      </w:r>
    </w:p>
  </w:tc>
</w:tr>
  <w:tr>
  <w:trPr>
    <w:trHeight w:val="933"/>
  </w:trPr>
  <w:tc>
    <w:tcPr>
      <w:tcW w:w="3310" w:type="dxa"/>
    </w:tcPr>
    <w:p>
      <w:pPr>
        <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>TB 1.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>8.6.1</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t xml:space="preserve">Statement  No. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>2</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:sz w:val="36"/>
          <w:szCs w:val="36"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t xml:space="preserve">Dasini  No. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
          <w:szCs w:val="32"/>
        </w:rPr>
        <w:t>21</w:t>
      </w:r>
    </w:p>
  </w:tc>
  <w:tc>
  <w:tcPr>
    <w:tcW w:w="4738" w:type="dxa"/>
  </w:tcPr>
  <w:p>
    <w:pPr>
      <w:widowControl w:val="0"/>
      <w:autoSpaceDE w:val="0"/>
      <w:autoSpaceDN w:val="0"/>
      <w:adjustRightInd w:val="0"/>
      <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
      <w:rPr>
        <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
        <w:sz w:val="36"/>
        <w:szCs w:val="36"/>
      </w:rPr>
    </w:pPr>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="Latha"/>
        <w:b/>
        <w:sz w:val="36"/>
        <w:szCs w:val="28"/>
      </w:rPr>
      <w:t>Æ</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:highlight w:val="yellow"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>வீ</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>ர்ய</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
        <w:b/>
        <w:sz w:val="36"/>
        <w:szCs w:val="28"/>
      </w:rPr>
      <w:t>†</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>ந்</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha" w:hint="cs"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>-</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>த</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:b/>
        <w:bCs/>
        <w:position w:val="-12"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>3</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>த்</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:b/>
        <w:bCs/>
        <w:position w:val="-12"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>4</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t xml:space="preserve">யாத் </w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
      </w:rPr>
      <w:t>|</w:t>
    </w:r>
  </w:p>
</w:tc>
  <w:tc>
  <w:tcPr>
    <w:tcW w:w="5220" w:type="dxa"/>
  </w:tcPr>
  <w:p>
    <w:pPr>
      <w:autoSpaceDE w:val="0"/>
      <w:autoSpaceDN w:val="0"/>
      <w:adjustRightInd w:val="0"/>
      <w:spacing w:before="0" w:line="252" w:lineRule="auto"/>
      <w:rPr>
        <w:rFonts w:cs="Arial"/>
        <w:sz w:val="36"/>
        <w:szCs w:val="36"/>
      </w:rPr>
    </w:pPr>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="Latha"/>
        <w:b/>
        <w:sz w:val="36"/>
        <w:szCs w:val="28"/>
      </w:rPr>
      <w:t>Æ</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:highlight w:val="green"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>வீ</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
        <w:b/>
        <w:sz w:val="32"/>
        <w:szCs w:val="28"/>
      </w:rPr>
      <w:t>…</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>ர்ய</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="BRH Tamil Tab Extra" w:hAnsi="BRH Tamil Tab Extra" w:cs="Latha"/>
        <w:b/>
        <w:sz w:val="36"/>
        <w:szCs w:val="28"/>
      </w:rPr>
      <w:t>†</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>ந்</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha" w:hint="cs"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>-</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>த</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:b/>
        <w:bCs/>
        <w:position w:val="-12"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>3</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>த்</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:b/>
        <w:bCs/>
        <w:position w:val="-12"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t>4</w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
        <w:cs/>
        <w:lang w:bidi="ta-IN"/>
      </w:rPr>
      <w:t xml:space="preserve">யாத் </w:t>
    </w:r>
    <w:r>
      <w:rPr>
        <w:rFonts w:ascii="Latha" w:hAnsi="Latha" w:cs="Latha"/>
        <w:sz w:val="28"/>
        <w:szCs w:val="28"/>
      </w:rPr>
      <w:t>|</w:t>
    </w:r>
  </w:p>
</w:tc>
</w:tr>
</w:tbl>
<w:p>
  <w:pPr>
    <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="32"/>
      <w:szCs w:val="32"/>
      <w:u w:val="single"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>'
$r = $d.Range(0, 0)
$r.InsertXML($xml)
